# feat: add 2022-Q4 data
#
# This workbook keeps one worksheet per reporting quarter, plus a "总计"
# (totals) summary sheet. We are adding a new "2022-Q4" quarter:
#   1. Duplicate the "2022-Q3" sheet (so formatting / layout matches the
#      other quarter sheets exactly), place the copy right before it,
#      rename it to "2022-Q4", and update its figures.
#   2. Insert a new row at the top of the "总计" summary sheet for the
#      2022-Q4 figures, and bump the running index in column A for every
#      row that shifted down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q4" quarter sheet by duplicating "2022-Q3"
# ---------------------------------------------------------------------
$wsQ3 = $wb.Worksheets.Item("2022-Q3")
$wsQ3.Copy($wsQ3)
$wsQ4 = $wb.Worksheets.Item("2022-Q3 (2)")
$wsQ4.Name = "2022-Q4"

# D:G hold numeric-looking figures that are stored as text in this
# workbook (e.g. "15.68"); force text format so the leading/trailing
# digits and trailing zeros are preserved instead of being coerced to
# a real number.
$wsQ4.Range("D2:G3").NumberFormat = "@"

# Row 2: fund 001643 (汇丰晋信智造先锋股票A)
$wsQ4.Range("D2").Value = "16.39"
$wsQ4.Range("E2").Value = "93.66"
$wsQ4.Range("F2").Value = "3.86"
$wsQ4.Range("G2").Value = "0.6327"
$wsQ4.Range("H2").Value = 7

# Row 3: fund 001644 (汇丰晋信智造先锋股票C)
$wsQ4.Range("D3").Value = "8.77"
$wsQ4.Range("E3").Value = "93.66"
$wsQ4.Range("F3").Value = "3.86"
$wsQ4.Range("G3").Value = "0.3385"
$wsQ4.Range("H3").Value = 7

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

# Insert a fresh row above the current first data row (row 2), pushing
# the existing quarters down by one.
$wsTotal.Rows.Item(2).Insert()

# The inserted row picks up stray formatting from the insert operation;
# clear it so the blank row starts out unstyled like the source did.
$wsTotal.Range("B2:D2").ClearFormats()

# Column A uses the same "index" style on every data row; copy it from
# the row that used to be first (now row 3) onto the new row 2.
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0.97

# Re-number the running index in column A for every row that shifted
# down (it must stay a contiguous 0-based sequence).
for ($r = 3; $r -le 9; $r++) {
    $wsTotal.Range("A$r").Value = $r - 2
}

Write-Output "2022-Q4 sheet added and 总计 summary updated"
